# Update "想去人数" (column F) counts across the four sheets of
# 北京-漫展信息.xlsx to reflect freshly scraped numbers.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws3 = $wb.Worksheets.Item(3)   # 本地生活
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# 展览 (sheet 1)
$ws1.Range("F5").Value = 71
$ws1.Range("F6").Value = 834
$ws1.Range("F7").Value = 418
$ws1.Range("F8").Value = 4691
$ws1.Range("F9").Value = 4691
$ws1.Range("F16").Value = 7466
$ws1.Range("F21").Value = 520
$ws1.Range("F22").Value = 1357
$ws1.Range("F24").Value = 6284
$ws1.Range("F25").Value = 1743
$ws1.Range("F28").Value = 6167
$ws1.Range("F34").Value = 6402
$ws1.Range("F45").Value = 39
$ws1.Range("F46").Value = 431
$ws1.Range("F47").Value = 2137
$ws1.Range("F49").Value = 1074

# 演出 (sheet 2)
$ws2.Range("F6").Value = 124

# 本地生活 (sheet 3)
$ws3.Range("F2").Value = 1444

# 全部类型 (sheet 4)
$ws4.Range("F3").Value = 1444
$ws4.Range("F6").Value = 71
$ws4.Range("F8").Value = 418
$ws4.Range("F9").Value = 4691
$ws4.Range("F10").Value = 4691
$ws4.Range("F17").Value = 7466
$ws4.Range("F20").Value = 520
$ws4.Range("F21").Value = 1357
$ws4.Range("F22").Value = 124
$ws4.Range("F23").Value = 6284
$ws4.Range("F24").Value = 1743
$ws4.Range("F29").Value = 6167
$ws4.Range("F36").Value = 6402
$ws4.Range("F45").Value = 39
$ws4.Range("F46").Value = 431
$ws4.Range("F48").Value = 2137
